# Update the "Förändrad" (last modified) date column for the logging
# report rows. All rows currently stamped 2023-10-13 (serial 45212) are
# refreshed to 2023-10-22 (serial 45221).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C2:C10").Value = 45221
